# Fruta / hortaliza, semanal
# Weekly update: two new price records added at the top of the data block
# (rows 15-16), pushing the previously-existing rows 15-23 down to 17-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 15, shifting the existing rows 15-23 down to 17-25.
$ws.Range("A15:A16").EntireRow.Insert()

# Columns A, B, C, E-K are constant for every record in this data block.
$commonA = 10
$commonB = "Vega Modelo de Temuco"
$commonC = "La Araucanía"
$commonE = 9
$commonF = "Fruta"
$commonG = 100104
$commonH = "Frutos de pepita"
$commonI = 100104004
$commonJ = "Níspero"
$commonK = "Californiana(o)"

# Final data for rows 15-25 (D, L, M, N, O, P, Q, R, S, T).
$rows = @(
  @{ R=15; D=45225; L="Primera"; M=200; N=35000; O=35000; P=35000; Q="$/bandeja 10 kilos"; Rg="Provincia de Los Andes"; S=3500;  T=10 },
  @{ R=16; D=45225; L="Segunda"; M=90;  N=24000; O=24000; P=24000; Q="$/bandeja 10 kilos"; Rg="Provincia de Los Andes"; S=2400;  T=10 },
  @{ R=17; D=44868; L="Primera"; M=30;  N=14000; O=14000; P=14000; Q="$/bandeja 5 kilos";  Rg="Provincia de Quillota";  S=2800;  T=5  },
  @{ R=18; D=45222; L="Primera"; M=25;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; Rg="Provincia de Quillota";  S=2800;  T=10 },
  @{ R=19; D=44889; L="Primera"; M=50;  N=30000; O=30000; P=30000; Q="$/bandeja 10 kilos"; Rg="Provincia de Quillota";  S=3000;  T=10 },
  @{ R=20; D=44511; L="Primera"; M=45;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; Rg="Provincia de Los Andes"; S=2800;  T=10 },
  @{ R=21; D=44511; L="Primera"; M=45;  N=3200;  O=3200;  P=3200;  Q="$/bandeja 10 kilos"; Rg="Provincia de Quillota";  S=320;   T=10 },
  @{ R=22; D=44503; L="Primera"; M=50;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; Rg="Provincia de Quillota";  S=2800;  T=10 },
  @{ R=23; D=44874; L="Primera"; M=40;  N=25000; O=25000; P=25000; Q="$/bandeja 10 kilos"; Rg="Provincia de Quillota";  S=2500;  T=10 },
  @{ R=24; D=44921; L="Primera"; M=55;  N=15000; O=15000; P=15000; Q="$/bandeja 7 kilos";  Rg="Provincia de Quillota";  S=2143;  T=7  },
  @{ R=25; D=45224; L="Primera"; M=40;  N=30000; O=30000; P=30000; Q="$/bandeja 10 kilos"; Rg="Provincia de Los Andes"; S=3000;  T=10 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $commonA
    $ws.Cells.Item($r, 2).Value = $commonB
    $ws.Cells.Item($r, 3).Value = $commonC
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $commonE
    $ws.Cells.Item($r, 6).Value = $commonF
    $ws.Cells.Item($r, 7).Value = $commonG
    $ws.Cells.Item($r, 8).Value = $commonH
    $ws.Cells.Item($r, 9).Value = $commonI
    $ws.Cells.Item($r, 10).Value = $commonJ
    $ws.Cells.Item($r, 11).Value = $commonK
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.Rg
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
